$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10; this shifts the existing rows 10-43 down
# to 11-44 (carrying all their data/formatting with them), matching the
# target diff where old row 10 becomes new row 11, ..., old row 43 becomes
# new row 44.
$ws.Rows.Item(10).Insert()

# Populate the newly-inserted row 10 with this week's new data record.
$ws.Cells.Item(10, 1).Value2 = 10
$ws.Cells.Item(10, 2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item(10, 3).Value2 = "La Araucanía"
$ws.Cells.Item(10, 4).Value2 = 44742
$ws.Cells.Item(10, 5).Value2 = 9
$ws.Cells.Item(10, 6).Value2 = "Fruta"
$ws.Cells.Item(10, 7).Value2 = 100108
$ws.Cells.Item(10, 8).Value2 = "Tropicales y subtropicales"
$ws.Cells.Item(10, 9).Value2 = 100108003
$ws.Cells.Item(10, 10).Value2 = "Maracuyá"
$ws.Cells.Item(10, 11).Value2 = "Sin especificar"
$ws.Cells.Item(10, 12).Value2 = "Primera"
$ws.Cells.Item(10, 13).Value2 = 5
$ws.Cells.Item(10, 14).Value2 = 34000
$ws.Cells.Item(10, 15).Value2 = 34000
$ws.Cells.Item(10, 16).Value2 = 34000
$ws.Cells.Item(10, 17).Value2 = "$/caja 18 kilos"
$ws.Cells.Item(10, 18).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(10, 19).Value2 = 1889
$ws.Cells.Item(10, 20).Value2 = 18
